$d = $word.ActiveDocument

# --- Title / Title Char: drop the explicit character spacing / kerning ---
# (Word's Font dialog "Normal" spacing + kerning off both correspond to the
# OOXML default value of 0 for <w:spacing> / <w:kern>.)
$title = $d.Styles("Title")
$title.Font.Spacing = 0
$title.Font.Kerning = 0

$titleChar = $d.Styles("Title Char")
$titleChar.Font.Spacing = 0
$titleChar.Font.Kerning = 0

# --- Author: now based on Title, no longer forces its own centering,
#     and carries an explicit 12pt (24 half-point) run size ---
$author = $d.Styles("Author")
$author.BaseStyle = $title
$author.Font.Size = 12
$author.Font.SizeBi = 12

# --- Date: same treatment as Author ---
$date = $d.Styles("Date")
$date.BaseStyle = $title
$date.Font.Size = 12
$date.Font.SizeBi = 12
